# Updates cryptos list values (Price / Volume(1h)) for rows 2-51.
# Column D values that parse as plain numbers are prefixed with a leading
# apostrophe so Excel stores them as text (matching the sheet's original
# text-typed cells) instead of silently converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.548.60"
$ws.Range("E2").Value = "  +5.69%  "
$ws.Range("D3").Value = "1.722.18"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'225.78"
$ws.Range("E5").Value = "  +3.52%  "
$ws.Range("D6").Value = "'0.5373"
$ws.Range("E6").Value = "  +3.07%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'0.2664"
$ws.Range("E8").Value = "  +1.12%  "
$ws.Range("D9").Value = "'0.06603"
$ws.Range("E9").Value = "  +4.35%  "
$ws.Range("D10").Value = "'21.66"
$ws.Range("D11").Value = "'0.07723"
$ws.Range("E11").Value = "  +0.72%  "
$ws.Range("D12").Value = "'4.615"
$ws.Range("E12").Value = "  +0.70%  "
$ws.Range("D13").Value = "1.710.92"
$ws.Range("E13").Value = "  +4.01%  "
$ws.Range("D14").Value = "1.960.13"
$ws.Range("E14").Value = "  +4.42%  "
$ws.Range("D15").Value = "'0.5839"
$ws.Range("E15").Value = "  +4.49%  "
$ws.Range("D16").Value = "0.0₅8318"
$ws.Range("E16").Value = "  +2.27%  "
$ws.Range("D17").Value = "'67.93"
$ws.Range("E17").Value = "  +4.12%  "
$ws.Range("D18").Value = "27.551.92"
$ws.Range("E18").Value = "  +5.74%  "
$ws.Range("D19").Value = "'220.39"
$ws.Range("E19").Value = "  +15.13%  "
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("D21").Value = "'4.724"
$ws.Range("E21").Value = "  +2.18%  "
$ws.Range("D22").Value = "'10.65"
$ws.Range("E22").Value = "  +1.59%  "
$ws.Range("D23").Value = "'6.079"
$ws.Range("E23").Value = "  +2.97%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").Value = "'148.30"
$ws.Range("E25").Value = "  +2.92%  "
$ws.Range("E26").Value = "  +14.93%  "
$ws.Range("D27").Value = "'0.1236"
$ws.Range("E27").Value = "  +4.19%  "
$ws.Range("D28").Value = "'7.405"
$ws.Range("E28").Value = "  +2.84%  "
$ws.Range("D29").Value = "'16.59"
$ws.Range("E29").Value = "  +4.45%  "
$ws.Range("D30").Value = "'0.05574"
$ws.Range("E30").Value = "  +2.56%  "
$ws.Range("E31").Value = "  +2.84%  "
$ws.Range("D32").Value = "'3.561"
$ws.Range("E32").Value = "  +3.50%  "
$ws.Range("D33").Value = "'3.442"
$ws.Range("E33").Value = "  +2.92%  "
$ws.Range("E34").Value = "  +6.98%  "
$ws.Range("E35").Value = "  +2.09%  "
$ws.Range("D36").Value = "'0.9636"
$ws.Range("E36").Value = "  +1.99%  "
$ws.Range("D37").Value = "'2.423"
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("D38").Value = "'0.5972"
$ws.Range("E38").Value = "  +6.03%  "
$ws.Range("E39").Value = "  +4.37%  "
$ws.Range("D40").Value = "'5.921"
$ws.Range("E40").Value = "  +0.87%  "
$ws.Range("D41").Value = "'0.8550"
$ws.Range("E41").Value = "  +3.39%  "
$ws.Range("D42").Value = "1.053.21"
$ws.Range("E42").Value = "  +2.62%  "
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("D44").Value = "'101.23"
$ws.Range("E44").Value = "  +0.26%  "
$ws.Range("D45").Value = "1.866.40"
$ws.Range("E45").Value = "  +4.55%  "
$ws.Range("E46").Value = "  +3.92%  "
$ws.Range("D47").Value = "'59.14"
$ws.Range("E47").Value = "  +3.18%  "
$ws.Range("D48").Value = "'8.245"
$ws.Range("E48").Value = "  +4.01%  "
$ws.Range("D50").Value = "'1.001"
$ws.Range("E50").Value = "  +0.25%  "
$ws.Range("D51").Value = "'0.05245"
$ws.Range("E51").Value = "  +2.29%  "
